$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values must remain plain text,
# matching the original inline-string cell type rather than becoming numbers/percentages).
$updates = [ordered]@{
    'D2' = '315.66'
    'E2' = '2.49%'
    'D3' = '39.39'
    'E3' = '-0.91%'
    'D4' = '5.138'
    'E4' = '-0.40%'
    'D5' = '0.08205'
    'E5' = '1.28%'
    'D6' = '1.995'
    'E6' = '2.66%'
    'D7' = '8.312'
    'E7' = '2.01%'
    'D8' = '0.9379'
    'E8' = '1.06%'
    'D9' = '0.1305'
    'E9' = '-8.85%'
    'D10' = '0.1974'
    'E10' = '2.72%'
    'D11' = '0.09034'
    'E11' = '-1.35%'
    'D12' = '0.03472'
    'E12' = '-1.02%'
    'E13' = '-0.63%'
    'D14' = '0.001416'
    'E14' = '0.94%'
    'D15' = '0.006498'
    'E15' = '10.48%'
    'E16' = '-7.71%'
    'D17' = '4.364'
    'E17' = '3.03%'
    'E18' = '-1.66%'
    'D19' = '0.3472'
    'E19' = '1.28%'
    'E20' = '-0.45%'
    'E21' = '6.33%'
    'D22' = '0.2488'
    'E22' = '2.67%'
    'D23' = '0.04366'
    'E23' = '-0.06%'
    'D24' = '0.001239'
    'E24' = '0.85%'
    'D25' = '0.004774'
    'E25' = '9.53%'
    'D26' = '0.0003888'
    'E26' = '198.68%'
    'E27' = '-7.62%'
    'D39' = '0.02236'
    'E39' = '9.73%'
    'D40' = '0.05214'
    'E40' = '2.89%'
    'D41' = '0.007743'
    'E41' = '5.20%'
    'D42' = '0.01035'
    'E42' = '5.96%'
    'E43' = '2.67%'
    'E44' = '-1.59%'
    'D45' = '0.009751'
    'E45' = '2.96%'
    'D46' = '0.00006746'
    'E46' = '6.30%'
    'E47' = '-0.12%'
    'D48' = '0.002881'
    'E48' = '5.56%'
    'E49' = '29.90%'
    'E50' = '-0.12%'
    'E51' = '-0.12%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage: set a text format first so Excel keeps the literal
    # string (e.g. "315.66" or "2.49%") instead of auto-converting it to a
    # number/percentage, then clear the temporary format so no new cell style
    # is left applied to the cell (matching the unstyled original cells).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
